# "Generate Report for Archive" - regenerate the localization-status report:
#   1. Flip the outstanding handoff rows from "Ready for handoff" to "In Translation"
#      (Overview sheet E2:F3, and the Status column on each per-locale sheet).
#   2. Shrink the "Status"-ish columns (E/F on Overview, C on the locale sheets) to
#      their freshly autofit width now that the cell text is shorter.

$wb = $excel.ActiveWorkbook

$oldStatus = "Ready for handoff"
$newStatus = "In Translation"

# --- Overview sheet -------------------------------------------------------
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E2").Value = $newStatus
$overview.Range("F2").Value = $newStatus
$overview.Range("E3").Value = $newStatus
$overview.Range("F3").Value = $newStatus

# Columns E and F narrow from ~17.22 chars to ~13.41 chars.
$overview.Columns.Item(5).ColumnWidth = 12.576851254417766
$overview.Columns.Item(6).ColumnWidth = 12.576851254417766

# --- Per-locale sheets (zh-cn, de-de): Status lives in column C -----------
foreach ($sheetName in @("zh-cn", "de-de")) {
    $ws = $wb.Worksheets.Item($sheetName)
    $ws.Range("C2").Value = $newStatus
    $ws.Range("C3").Value = $newStatus

    # Column C narrows from ~17.22 chars to ~13.41 chars.
    $ws.Columns.Item(3).ColumnWidth = 12.576851254417766
}
